$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3, shifting the existing row 3 (B3) and
# row 4 (B4,D4,E4) down to rows 4 and 5 respectively.
$null = $ws.Rows(3).Insert()

# Give the now-empty A1 cell a thin box border on all four sides. This
# creates the formatted-but-blank cell (new style xf with borderId=1)
# used to test that empty-but-formatted cells aren't treated as blank.
$ws.Range("A1").Borders.LineStyle = 1

# Match the author's final selection/active cell.
$null = $ws.Range("A4:E5").Select()
